$wb = $excel.ActiveWorkbook

# --- comp_quantity_inst1 sheet: add "type" column (E) with letters A..O ---
$ws3 = $wb.Worksheets.Item("comp_quantity_inst1")

$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $ws3.Range("E$row").Value = $letters[$i]
}
$ws3.Range("E1").Value = "type"

# Center-align the whole used range A1:E16 like the new data layout
$ws3.Range("A1:E16").HorizontalAlignment = -4108

# --- parameters sheet: update totals ---
$ws1 = $wb.Worksheets.Item("parameters")
$ws1.Range("B12").Value = 25

$ws1.Range("A13").Value = "no_req_total"
$ws1.Range("B13").Formula = "=SUM(comp_quantity_inst1!C2:C9)"

$ws1.Range("A14").Value = "no_opt_total"
$ws1.Range("B14").Formula = "=SUM(comp_quantity_inst1!C10:C16)"

# --- selection / active sheet state ---
$ws3.Range("C19").Select()

$ws1.Activate()
$ws1.Range("A13:A14").Select()
